$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1 / Q1, copying the style of O1 (bold, centered, bordered) ---
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25 ---
# For each row: swap I<->K and M<->O values, and append P=2, Q=2 (copying N's plain style)
for ($r = 2; $r -le 25; $r++) {
    $iCell = $ws.Cells.Item($r, 9)   # I
    $kCell = $ws.Cells.Item($r, 11)  # K
    $mCell = $ws.Cells.Item($r, 13)  # M
    $oCell = $ws.Cells.Item($r, 15)  # O

    $iVal = $iCell.Formula
    $kVal = $kCell.Formula
    $mVal = $mCell.Formula
    $oVal = $oCell.Formula

    $iCell.Formula = $kVal
    $kCell.Formula = $mVal
    $mCell.Formula = $oVal
    $oCell.Formula = $iVal

    $nAddr = "N" + $r
    $pAddr = "P" + $r
    $qAddr = "Q" + $r

    $ws.Range($nAddr).Copy($ws.Range($pAddr))
    $ws.Range($nAddr).Copy($ws.Range($qAddr))

    $ws.Range($pAddr).Value = 2
    $ws.Range($qAddr).Value = 2
}
